$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 6816.1665
$ws.Range("I43").Value = 6479.4
$ws.Range("J43").Value = 8500
$ws.Range("K43").Value = 6479.4
$ws.Range("L43").Value = 8500
$ws.Range("M43").Value = -6410.4
$ws.Range("N43").Value = -8638

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 47320.94
$ws.Range("J86").Value = 116292.4
$ws.Range("L86").Value = 116292.4
$ws.Range("N86").Value = -118538.4

# Row 88: The Grave of Hemlock Groves / Growth Formula Zeta
$ws.Range("H88").Value = 41017268
$ws.Range("I88").Value = 111111784
$ws.Range("J88").Value = 5970013.5
$ws.Range("K88").Value = 111111784
$ws.Range("L88").Value = 5970013.5
$ws.Range("M88").Value = -111111378
$ws.Range("N88").Value = -5970825.5

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 47320.94
$ws.Range("J89").Value = 116292.4
$ws.Range("L89").Value = 581462
$ws.Range("N89").Value = -592694

# Row 91: Dappling the Highlands (L) / Growth Formula Zeta
$ws.Range("H91").Value = 41017268
$ws.Range("I91").Value = 111111784
$ws.Range("J91").Value = 5970013.5
$ws.Range("K91").Value = 111111784
$ws.Range("L91").Value = 5970013.5
$ws.Range("M91").Value = -111110380
$ws.Range("N91").Value = -5972821.5

# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws.Range("H92").Value = 799.63635
$ws.Range("I92").Value = 950.7778
$ws.Range("J92").Value = 119.5
$ws.Range("K92").Value = 950.7778
$ws.Range("L92").Value = 119.5
$ws.Range("M92").Value = 297.2222
$ws.Range("N92").Value = -2615.5

# Row 93: Spellbound / Koppranickel Index
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents() | Out-Null

# Row 97: Materia Worth / Potent Spiritbond Potion
$ws.Range("H97").Value = 4382
$ws.Range("J97").Value = 4382
$ws.Range("L97").Value = 13146
$ws.Range("N97").Value = -14138

# Row 123: Nearly Bare / Gaja Grimoire
$ws.Range("H123").Value = 67009.164
$ws.Range("J123").Value = 67009.164
$ws.Range("L123").Value = 67009.164
$ws.Range("N123").Value = -76809.164

# Row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 1200
$ws.Range("I131").Value = 1200
$ws.Range("K131").Value = 3600
$ws.Range("M131").Value = 1440

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2674.5366
$ws.Range("I138").Value = 2421.2942
$ws.Range("K138").Value = 7263.882599999999
$ws.Range("M138").Value = -2123.882599999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 4455.7144
$ws.Range("I32").Value = 2898.5
$ws.Range("J32").Value = 11376.667
$ws.Range("K32").Value = 2898.5
$ws.Range("L32").Value = 11376.667
$ws.Range("M32").Value = -2611.5
$ws.Range("N32").Value = -11950.667

# Row 95: Shielded Life / High Steel Scutum
$ws.Range("H95").Value = 60104
$ws.Range("J95").Value = 60104
$ws.Range("L95").Value = 60104
$ws.Range("N95").Value = -65596

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall / Elm Lumber
$ws.Range("H22").Value = 25349.5
$ws.Range("I22").Value = 50199.5
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 50199.5
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = -49849.5
$ws.Range("N22").Value = -1199.5

# Row 75: The Darkest Hearth / Dark Chestnut Spinning Wheel
$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41996

# Row 78: Fruit of the Loom (L) / Dark Chestnut Spinning Wheel
$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129984

# Row 107: Built to Last / White Oak Lumber
$ws.Range("I107").Value = 362630.03
$ws.Range("J107").Value = 64763.938
$ws.Range("K107").Value = 362630.03
$ws.Range("L107").Value = 64763.938
$ws.Range("M107").Value = -360710.03
$ws.Range("N107").Value = -68603.93799999999

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 25001354
$ws.Range("I132").Value = 25001354
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 75004062
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -75001532
$ws.Range("N132").ClearContents() | Out-Null

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 11365233
$ws.Range("I134").Value = 13890224
$ws.Range("K134").Value = 41670672
$ws.Range("M134").Value = -41668137

# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 445938.4
$ws.Range("J141").Value = 532423
$ws.Range("L141").Value = 532423
$ws.Range("N141").Value = -542783

$ws = $wb.Worksheets.Item("CUL")
# Row 37: I Love Lamprey / Eel Pie
$ws.Range("H37").Value = 137868
$ws.Range("J37").Value = 137868
$ws.Range("L37").Value = 413604
$ws.Range("N37").Value = -413828

# Row 112: Sweet Tooth / Caramels
$ws.Range("H112").Value = 14480
$ws.Range("I112").Value = 7705
$ws.Range("J112").Value = 19900
$ws.Range("K112").Value = 23115
$ws.Range("L112").Value = 59700
$ws.Range("M112").Value = -22007
$ws.Range("N112").Value = -61916

# Row 128: A Historical Flavor / Skyr
$ws.Range("H128").Value = 138921.28
$ws.Range("I128").Value = 138921.28
$ws.Range("K128").Value = 416763.84
$ws.Range("M128").Value = -411783.84

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1777
$ws.Range("I131").Value = 1740
$ws.Range("K131").Value = 5220
$ws.Range("M131").Value = -180

$ws = $wb.Worksheets.Item("GSM")
# Row 58: The Big Red / Red Coral Necklace
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents() | Out-Null

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3131.879
$ws.Range("I126").Value = 3271.138
$ws.Range("J126").Value = 2122.25
$ws.Range("K126").Value = 9813.414000000001
$ws.Range("L126").Value = 6366.75
$ws.Range("M126").Value = -7343.414000000001
$ws.Range("N126").Value = -11306.75

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 6584050
$ws.Range("I132").Value = 7356696.5
$ws.Range("K132").Value = 22070089.5
$ws.Range("M132").Value = -22067559.5

# Row 139: Ringing Gratitude / White Gold Ring of Healing
$ws.Range("H139").Value = 112500
$ws.Range("J139").Value = 112500
$ws.Range("L139").Value = 112500
$ws.Range("N139").Value = -122780

$ws = $wb.Worksheets.Item("LTW")
# Row 6: Sticking Their Necks Out / Leather Choker
$ws.Range("H6").Value = 54685.09
$ws.Range("I6").Value = 44995
$ws.Range("J6").Value = 55654.1
$ws.Range("K6").Value = 44995
$ws.Range("L6").Value = 55654.1
$ws.Range("M6").Value = -44883
$ws.Range("N6").Value = -55878.1

# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 4298.8
$ws.Range("I7").Value = 3776.4443
$ws.Range("K7").Value = 3776.4443
$ws.Range("M7").Value = -3664.4443

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 2483
$ws.Range("I22").Value = 2386.3333
$ws.Range("J22").Value = 2591.75
$ws.Range("K22").Value = 2386.3333
$ws.Range("L22").Value = 2591.75
$ws.Range("M22").Value = -2091.3333
$ws.Range("N22").Value = -3181.75

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 2483
$ws.Range("I27").Value = 2386.3333
$ws.Range("J27").Value = 2591.75
$ws.Range("K27").Value = 2386.3333
$ws.Range("L27").Value = 2591.75
$ws.Range("M27").Value = -2279.3333
$ws.Range("N27").Value = -2805.75

# Row 31: Open to Attack / Goatskin Jacket
$ws.Range("H31").Value = 6369.4
$ws.Range("J31").Value = 14333
$ws.Range("L31").Value = 14333
$ws.Range("N31").Value = -14829

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 900.35297
$ws.Range("I46").Value = 831.625
$ws.Range("K46").Value = 831.625
$ws.Range("M46").Value = -643.625

# Row 93: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 11743.235
$ws.Range("I93").Value = 13062.667
$ws.Range("J93").Value = 1847.5
$ws.Range("K93").Value = 13062.667
$ws.Range("L93").Value = 1847.5
$ws.Range("M93").Value = -11814.667
$ws.Range("N93").Value = -4343.5

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 4298.8
$ws.Range("I126").Value = 3776.4443
$ws.Range("K126").Value = 11329.3329
$ws.Range("M126").Value = -8859.332900000001

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 16560895
$ws.Range("I132").Value = 17787488
$ws.Range("K132").Value = 53362464
$ws.Range("M132").Value = -53359934

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 11631994
$ws.Range("I132").Value = 16130561
$ws.Range("J132").Value = 10696.417
$ws.Range("K132").Value = 48391683
$ws.Range("L132").Value = 32089.251
$ws.Range("M132").Value = -48389153
$ws.Range("N132").Value = -37149.251
